$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.067.59"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3
$ws.Range("D3").Value = "1.846.51"
$ws.Range("E3").Value = "  +2.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'233.52"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "  +2.77%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'41.85"
$ws.Range("E8").Value = "  +6.47%  "

# Row 9
$ws.Range("D9").Value = "'0.329"
$ws.Range("E9").Value = "  +0.91%  "

# Row 10
$ws.Range("D10").Value = "'0.0696"
$ws.Range("E10").Value = "  +2.05%  "

# Row 11
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "  -1.02%  "

# Row 12
$ws.Range("D12").Value = "2.112.78"
$ws.Range("E12").Value = "  +2.05%  "

# Row 13
$ws.Range("D13").Value = "'11.52"
$ws.Range("E13").Value = "  +3.94%  "

# Row 14
$ws.Range("D14").Value = "1.843.26"
$ws.Range("E14").Value = "  +1.67%  "

# Row 15
$ws.Range("D15").Value = "'0.677"
$ws.Range("E15").Value = "  +0.48%  "

# Row 16
$ws.Range("E16").Value = "  +2.55%  "

# Row 17
$ws.Range("D17").Value = "35.057.37"
$ws.Range("E17").Value = "  +1.07%  "

# Row 18
$ws.Range("D18").Value = "'70.01"
$ws.Range("E18").Value = "  +0.80%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").Value = "  +0.83%  "

# Row 20
$ws.Range("D20").Value = "'240.78"
$ws.Range("E20").Value = "  +0.32%  "

# Row 21
$ws.Range("D21").Value = "'12.18"
$ws.Range("E21").Value = "  +1.80%  "

# Row 22
$ws.Range("E22").Value = "  +2.33%  "

# Row 23
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("E24").Value = "  +3.30%  "

# Row 25
$ws.Range("D25").Value = "'171.99"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").Value = "'7.93"
$ws.Range("E26").Value = "  +2.85%  "

# Row 27
$ws.Range("D27").Value = "'17.54"
$ws.Range("E27").Value = "  +2.03%  "

# Row 28
$ws.Range("E28").Value = "  +3.79%  "

# Row 29
$ws.Range("E29").Value = "  +11.67%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0557"
$ws.Range("E30").Value = "  +2.22%  "

# Row 31
$ws.Range("B31").Value = "BinanceUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D31").Value = "'1.01"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32
$ws.Range("D32").Value = "'3.98"
$ws.Range("E32").Value = "  -1.20%  "

# Row 33
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("E34").Value = "  +23.83%  "

# Row 35
$ws.Range("D35").Value = "'1.98"
$ws.Range("E35").Value = "  +11.28%  "

# Row 36
$ws.Range("D36").Value = "'0.763"
$ws.Range("E36").Value = "  +9.20%  "

# Row 37
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "  -3.88%  "

# Row 38
$ws.Range("E38").Value = "  +11.62%  "

# Row 39
$ws.Range("D39").Value = "'90.55"
$ws.Range("E39").Value = "  -1.04%  "

# Row 40
$ws.Range("E40").Value = "  +4.51%  "

# Row 41
$ws.Range("D41").Value = "1.346.76"
$ws.Range("E41").Value = "  +1.81%  "

# Row 42
$ws.Range("D42").Value = "'14.64"
$ws.Range("E42").Value = "  +3.43%  "

# Row 43
$ws.Range("D43").Value = "'2.30"
$ws.Range("E43").Value = "  +4.05%  "

# Row 44
$ws.Range("E44").Value = "  +3.42%  "

# Row 45
$ws.Range("E45").Value = "  -3.82%  "

# Row 46
$ws.Range("D46").Value = "'0.0533"
$ws.Range("E46").Value = "  +3.97%  "

# Row 47
$ws.Range("D47").Value = "'6.35"
$ws.Range("E47").Value = "  +1.22%  "

# Row 48
$ws.Range("D48").Value = "'11.70"
$ws.Range("E48").Value = "  +70.15%  "

# Row 49
$ws.Range("D49").Value = "2.027.32"
$ws.Range("E49").Value = "  +1.51%  "

# Row 50
$ws.Range("D50").Value = "'3.45"
$ws.Range("E50").Value = "  +16.33%  "

# Row 51
$ws.Range("D51").Value = "'0.0673"
$ws.Range("E51").Value = "  +0.79%  "
